$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4
$ws.Range("C2").Value = "0.79 [0.57–0.95]"
$ws.Range("D2").Value = "66% [33%–100%]"
$ws.Range("E2").Value = "74% [62%–82%]"

$ws.Range("B20").Value = 0.02
$ws.Range("C20").Value = "0.79 [0.57–0.95]"
$ws.Range("D20").Value = "2% [0%–20%]"
$ws.Range("E20").Value = "100% [99%–100%]"
